$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, reusing the formatting of the existing header
# cells (e.g. G1) by copying the cell (value+format) then overwriting the
# value with the new header text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Fill H2:H6 with 0 (plain numeric values, default/no style - like B2:G6)
$ws.Range("H2:H6").Value = 0
